$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (everything currently in A:F shifts to B:G)
$ws.Columns("A:A").Insert()

# Copy the formatting (fill/font/style) of the header that just shifted into B1
# onto the new A1 cell, then set its text to "ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "ID"

# Give the new ID column its own width (~28.43 characters, matching the
# author's manual resize; ColumnWidth expects "screen" character units
# which the host quantizes to whole pixels on write-back)
$ws.Columns("A:A").ColumnWidth = 28.42578125 - 5/6

# Match the author's final selection
$ws.Range("B6").Select()
